# Add a new "2022-Q3" sheet right after "总计" (i.e. before the current
# "2022-Q2" sheet, which is the 2nd sheet) and populate it with the
# quarterly fund-holdings breakdown, then update the "总计" summary sheet
# with the new row + shifted history.

$wb = $excel.ActiveWorkbook

# --- 1. Insert the new sheet -------------------------------------------------
# NOTE: `Worksheets.Item(N)` resolves lazily against whatever sheet
# currently sits at position N - it is NOT a stable handle to a particular
# sheet. So finish every sheet-add/rename/move first, and only *afterwards*
# fetch the references we'll actually use for the per-cell work below.
$q3Sheet = $wb.Worksheets.Add($wb.Worksheets.Item(2))
$q3Sheet.Name = "2022-Q3"

$q3Sheet = $wb.Worksheets.Item(2)   # re-fetch: now stable, no more sheet-level ops follow
$q2Sheet = $wb.Worksheets.Item(3)   # "2022-Q2", pushed one slot down by the insert above

# Bring over the header/index-column formatting from the neighbouring
# "2022-Q2" sheet template (a brand new sheet starts out unstyled).
$q2Sheet.Range("B1:H1").Copy()
$q3Sheet.Range("B1:H1").PasteSpecial(-4122)
$q2Sheet.Range("A2:A15").Copy()
$q3Sheet.Range("A2:A15").PasteSpecial(-4122)

# --- 2. Fill in the 2022-Q3 fund holdings table -----------------------------
# Header text never looks numeric, so no Text-formatting is needed for row 1
# (doing so would also strip the bold/border header style copied above).
$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($c = 0; $c -lt $headers.Length; $c++) {
    $q3Sheet.Cells.Item(1, $c + 2).Value = $headers[$c]
}

$rows = @(
    @(0, "008269", "大成睿享混合A",                         "19.80", "66.91", "3.97", "0.7861", 5),
    @(1, "090013", "大成竞争优势混合",                       "6.88",  "61.00", "3.47", "0.2387", 5),
    @(2, "008270", "大成睿享混合C",                         "4.02",  "66.91", "3.97", "0.1596", 5),
    @(3, "013463", "大成致远优势一年持有期混合A",             "3.65",  "60.88", "2.81", "0.1026", 9),
    @(4, "011834", "大成投资严选六月持有混合A",               "3.10",  "66.75", "2.86", "0.0887", 8),
    @(5, "004317", "前海开源沪港深裕鑫灵活配置混合C",         "1.79",  "70.17", "2.40", "0.0430", 8),
    @(6, "004316", "前海开源沪港深裕鑫灵活配置混合A",         "1.77",  "70.17", "2.40", "0.0425", 8),
    @(7, "004098", "前海开源港股通股息率50强股票",            "0.31",  "87.14", "3.12", "0.0097", 7),
    @(8, "011835", "大成投资严选六月持有混合C",               "0.22",  "66.75", "2.86", "0.0063", 8),
    @(9, "013464", "大成致远优势一年持有期混合C",             "0.17",  "60.88", "2.81", "0.0048", 9),
    @(10, "501303", "广发恒生中型股指数（LOF）A",             "0.21",  "89.12", "1.32", "0.0028", 10),
    @(11, "004996", "广发恒生中型股指数（LOF）C",             "0.09",  "89.12", "1.32", "0.0012", 10),
    @(12, "160922", "大成恒生综合中小型股指数（QDII-LOF）A",  "0.09",  "86.62", "0.94", "0.0008", 10),
    @(13, "008972", "大成恒生综合中小型股指数C",              "0.02",  "86.62", "0.94", "0.0002", 10)
)

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $i + 2
    $row = $rows[$i]
    $q3Sheet.Range("B" + $r + ":G" + $r).NumberFormat = "@"
    $q3Sheet.Cells.Item($r, 1).Value = $row[0]
    $q3Sheet.Cells.Item($r, 2).Value = $row[1]
    $q3Sheet.Cells.Item($r, 3).Value = $row[2]
    $q3Sheet.Cells.Item($r, 4).Value = $row[3]
    $q3Sheet.Cells.Item($r, 5).Value = $row[4]
    $q3Sheet.Cells.Item($r, 6).Value = $row[5]
    $q3Sheet.Cells.Item($r, 7).Value = $row[6]
    $q3Sheet.Cells.Item($r, 8).Value = $row[7]
}

# --- 3. Update the "总计" summary sheet: insert a new leading row for
#        2022-Q3 and push the older quarters down one row. --------------------
$summary = $wb.Worksheets.Item(1)
$summary.Rows.Item(2).Insert()

# The row-above format copy (Excel's default Insert behaviour) only
# touched B2:D2 (since A1 is blank in the header row) - strip that and
# bring in A2's style explicitly from the row below so the "index" column
# keeps its usual look.
$summary.Range("B2:D2").ClearFormats()
$summary.Cells.Item(3, 1).Copy()
$summary.Cells.Item(2, 1).PasteSpecial(-4122)

$summary.Cells.Item(2, 1).Value = 0
$summary.Cells.Item(2, 2).Value = "2022-Q3"
$summary.Cells.Item(2, 3).Value = 14
$summary.Cells.Item(2, 4).Value = 1.49

for ($r = 3; $r -le 8; $r++) {
    $summary.Cells.Item($r, 1).Value = $r - 2
}
